$d = $word.ActiveDocument
Write-Host "Test"
